$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 19999
$ws.Range("J12").Value = 19999
$ws.Range("L12").Value = 19999
$ws.Range("N12").Value = -20345

$ws.Range("H32").Value = 8023.194
$ws.Range("I32").Value = 5129.6577
$ws.Range("J32").Value = 11814.725
$ws.Range("K32").Value = 5129.6577
$ws.Range("L32").Value = 11814.725
$ws.Range("M32").Value = -4842.6577
$ws.Range("N32").Value = -12388.725

$ws.Range("H45").Value = 1468.1875
$ws.Range("I45").Value = 880.1177
$ws.Range("J45").Value = 2134.6667
$ws.Range("K45").Value = 880.1177
$ws.Range("L45").Value = 2134.6667
$ws.Range("M45").Value = -503.1177
$ws.Range("N45").Value = -2888.6667

$ws.Range("H61").Value = 2354.158
$ws.Range("I61").Value = 1292.3334
$ws.Range("J61").Value = 3309.8
$ws.Range("K61").Value = 1292.3334
$ws.Range("L61").Value = 3309.8
$ws.Range("M61").Value = -1080.3334
$ws.Range("N61").Value = -3733.8

$ws.Range("H74").Value = 3448.3684
$ws.Range("I74").Value = 3364.8386
$ws.Range("K74").Value = 3364.8386
$ws.Range("M74").Value = -2490.8386

$ws.Range("H77").Value = 3448.3684
$ws.Range("I77").Value = 3364.8386
$ws.Range("K77").Value = 16824.193
$ws.Range("M77").Value = -12456.193

$ws.Range("H132").Value = 3414.0322
$ws.Range("I132").Value = 2589.3809
$ws.Range("J132").Value = 5145.8
$ws.Range("K132").Value = 7768.1427
$ws.Range("L132").Value = 15437.4
$ws.Range("M132").Value = -5238.1427
$ws.Range("N132").Value = -20497.4

$ws.Range("H136").Value = 2354.158
$ws.Range("I136").Value = 1292.3334
$ws.Range("J136").Value = 3309.8
$ws.Range("K136").Value = 3877.0002
$ws.Range("L136").Value = 9929.400000000001
$ws.Range("M136").Value = -1327.0002
$ws.Range("N136").Value = -15029.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2978.5186
$ws.Range("I99").Value = 1320
$ws.Range("K99").Value = 1320
$ws.Range("M99").Value = 178

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2462.75
$ws.Range("I69").Value = 1373
$ws.Range("J69").Value = 3552.5
$ws.Range("K69").Value = 4119
$ws.Range("L69").Value = 10657.5
$ws.Range("M69").Value = -3308
$ws.Range("N69").Value = -12279.5

$ws.Range("H72").Value = 2462.75
$ws.Range("I72").Value = 1373
$ws.Range("J72").Value = 3552.5
$ws.Range("K72").Value = 12357
$ws.Range("L72").Value = 31972.5
$ws.Range("M72").Value = -8301
$ws.Range("N72").Value = -40084.5

$ws.Range("H80").Value = 6574.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 6574.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 19723.5
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -21595.5

$ws.Range("H83").Value = 6574.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 6574.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 59170.5
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -68530.5

$ws.Range("H131").Value = 10870564
$ws.Range("I131").Value = 125002720
$ws.Range("J131").Value = 835.0714
$ws.Range("K131").Value = 375008160
$ws.Range("L131").Value = 2505.2142
$ws.Range("M131").Value = -375003120
$ws.Range("N131").Value = -12585.2142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6830.3335
$ws.Range("I70").Value = 6030.95
$ws.Range("K70").Value = 6030.95
$ws.Range("M70").Value = -5760.95

$ws.Range("H73").Value = 6830.3335
$ws.Range("I73").Value = 6030.95
$ws.Range("K73").Value = 6030.95
$ws.Range("M73").Value = -5094.95

$ws.Range("H137").Value = 72819.78
$ws.Range("J137").Value = 72819.78
$ws.Range("L137").Value = 72819.78
$ws.Range("N137").Value = -83019.78

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1975.0533
$ws.Range("J2").Value = 1975.0533
$ws.Range("L2").Value = 1975.0533
$ws.Range("N2").Value = -2199.0533

$ws.Range("H61").Value = 1944.7059
$ws.Range("I61").Value = 1904.5454
$ws.Range("J61").Value = 2018.3334
$ws.Range("K61").Value = 1904.5454
$ws.Range("L61").Value = 2018.3334
$ws.Range("M61").Value = -1702.5454
$ws.Range("N61").Value = -2422.3334

$ws.Range("H113").Value = 1944.7059
$ws.Range("I113").Value = 1904.5454
$ws.Range("J113").Value = 2018.3334
$ws.Range("K113").Value = 1904.5454
$ws.Range("L113").Value = 2018.3334
$ws.Range("M113").Value = 265.4546
$ws.Range("N113").Value = -6358.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119:N133").ClearContents()
$ws.Range("H135:N141").ClearContents()
